$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the "25M" label (used in D1 and J1) with "25F".
# Excel will drop the now-unused "25M" shared string and append
# the new "25F" string, shifting the other shared-string indices
# (Strobe / Key under local) down automatically.
$ws.Range("D1").Value = "25F"
$ws.Range("J1").Value = "25F"

# Move/restore the active selection to J1 (was B13).
[void]$ws.Range("J1").Select()
